$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add new row 18 with date, time, file name and observation
# Copy formatting from row 17 (A17:B17) down to row 18 so the date/time
# number formats (styles 1 and 3) are reused instead of creating new ones.
$ws.Range("A17:B17").Copy()
$ws.Range("A18:B18").PasteSpecial(-4122)  # xlPasteFormats

$ws.Cells.Item(18, 1).Value = 45756
$ws.Cells.Item(18, 2).Value = 0.63472222222222219

$ws.Cells.Item(18, 3).Value = "Futconnect0904 1514"
$ws.Cells.Item(18, 4).Value = "Manual, participação efetiva e filtros nas abas de performance."

# Update selection to reflect the new active cell after the edit
$ws.Range("D19").Select()
